$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats / styles) of the last existing data row
# (row 774) down across the new rows so the new cells inherit the date /
# number cell styles used throughout the table.
$ws.Range("A774:C774").Copy() | Out-Null
$ws.Range("A775:C802").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new rows (775-802) with the updated daily UF / IVP values.
$ws.Cells.Item(775, 1).Value = 44237
$ws.Cells.Item(775, 2).Value = 29156.34
$ws.Cells.Item(775, 3).Value = 30342.38
$ws.Cells.Item(776, 1).Value = 44238
$ws.Cells.Item(776, 2).Value = 29163.61
$ws.Cells.Item(776, 3).Value = 30346.26
$ws.Cells.Item(777, 1).Value = 44239
$ws.Cells.Item(777, 2).Value = 29170.87
$ws.Cells.Item(777, 3).Value = 30350.14
$ws.Cells.Item(778, 1).Value = 44240
$ws.Cells.Item(778, 2).Value = 29178.14
$ws.Cells.Item(778, 3).Value = 30354.02
$ws.Cells.Item(779, 1).Value = 44241
$ws.Cells.Item(779, 2).Value = 29185.41
$ws.Cells.Item(779, 3).Value = 30357.91
$ws.Cells.Item(780, 1).Value = 44242
$ws.Cells.Item(780, 2).Value = 29192.68
$ws.Cells.Item(780, 3).Value = 30361.79
$ws.Cells.Item(781, 1).Value = 44243
$ws.Cells.Item(781, 2).Value = 29199.96
$ws.Cells.Item(781, 3).Value = 30365.67
$ws.Cells.Item(782, 1).Value = 44244
$ws.Cells.Item(782, 2).Value = 29207.23
$ws.Cells.Item(782, 3).Value = 30369.56
$ws.Cells.Item(783, 1).Value = 44245
$ws.Cells.Item(783, 2).Value = 29214.51
$ws.Cells.Item(783, 3).Value = 30373.439999999999
$ws.Cells.Item(784, 1).Value = 44246
$ws.Cells.Item(784, 2).Value = 29221.79
$ws.Cells.Item(784, 3).Value = 30377.33
$ws.Cells.Item(785, 1).Value = 44247
$ws.Cells.Item(785, 2).Value = 29229.07
$ws.Cells.Item(785, 3).Value = 30381.21
$ws.Cells.Item(786, 1).Value = 44248
$ws.Cells.Item(786, 2).Value = 29236.35
$ws.Cells.Item(786, 3).Value = 30385.1
$ws.Cells.Item(787, 1).Value = 44249
$ws.Cells.Item(787, 2).Value = 29243.64
$ws.Cells.Item(787, 3).Value = 30388.98
$ws.Cells.Item(788, 1).Value = 44250
$ws.Cells.Item(788, 2).Value = 29250.92
$ws.Cells.Item(788, 3).Value = 30392.87
$ws.Cells.Item(789, 1).Value = 44251
$ws.Cells.Item(789, 2).Value = 29258.21
$ws.Cells.Item(789, 3).Value = 30396.76
$ws.Cells.Item(790, 1).Value = 44252
$ws.Cells.Item(790, 2).Value = 29265.5
$ws.Cells.Item(790, 3).Value = 30400.65
$ws.Cells.Item(791, 1).Value = 44253
$ws.Cells.Item(791, 2).Value = 29272.79
$ws.Cells.Item(791, 3).Value = 30404.53
$ws.Cells.Item(792, 1).Value = 44254
$ws.Cells.Item(792, 2).Value = 29280.09
$ws.Cells.Item(792, 3).Value = 30408.42
$ws.Cells.Item(793, 1).Value = 44255
$ws.Cells.Item(793, 2).Value = 29287.38
$ws.Cells.Item(793, 3).Value = 30412.31
$ws.Cells.Item(794, 1).Value = 44256
$ws.Cells.Item(794, 2).Value = 29294.68
$ws.Cells.Item(794, 3).Value = 30416.2
$ws.Cells.Item(795, 1).Value = 44257
$ws.Cells.Item(795, 2).Value = 29301.98
$ws.Cells.Item(795, 3).Value = 30420.09
$ws.Cells.Item(796, 1).Value = 44258
$ws.Cells.Item(796, 2).Value = 29309.279999999999
$ws.Cells.Item(796, 3).Value = 30423.98
$ws.Cells.Item(797, 1).Value = 44259
$ws.Cells.Item(797, 2).Value = 29316.58
$ws.Cells.Item(797, 3).Value = 30427.87
$ws.Cells.Item(798, 1).Value = 44260
$ws.Cells.Item(798, 2).Value = 29323.89
$ws.Cells.Item(798, 3).Value = 30431.77
$ws.Cells.Item(799, 1).Value = 44261
$ws.Cells.Item(799, 2).Value = 29331.19
$ws.Cells.Item(799, 3).Value = 30435.66
$ws.Cells.Item(800, 1).Value = 44262
$ws.Cells.Item(800, 2).Value = 29338.5
$ws.Cells.Item(800, 3).Value = 30439.55
$ws.Cells.Item(801, 1).Value = 44263
$ws.Cells.Item(801, 2).Value = 29345.81
$ws.Cells.Item(801, 3).Value = 30443.45
$ws.Cells.Item(802, 1).Value = 44264
$ws.Cells.Item(802, 2).Value = 29353.119999999999
$ws.Cells.Item(802, 3).Value = 30447.34

# Update the named range so it covers the newly added rows.
$wb.Names.Item("UF_IVP_DIARIO").RefersTo = "=UF_IVP_DIARIO!`$A`$1:`$C`$802"

# Move the selection on the frozen (bottom-right) pane to the new last cell,
# matching where Excel leaves the cursor after appending data.
$ws.Range("A802").Select() | Out-Null
